# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# This updates column G ("K") for rows 2-39 with the newly-regenerated values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @(0,7,5,11,3,7,12,9,8,2,5,7,5,8,9,2,5,9,8,4,9,6,4,10,7,9,8,13,7,9,11,5,7,5,6,3,3,1)

$startRow = 2
for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Range("G$row").Value = $kValues[$i]
}
